$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.447.92'
$ws.Range("E2").Value = '  -3.46%  '
$ws.Range("D3").Value = '1.802.36'
$ws.Range("E3").Value = '  -3.20%  '
$ws.Range("D4").Value = '''1.009'
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''307.73'
$ws.Range("E6").Value = '  -2.49%  '
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("D8").Value = '''0.3645'
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").Value = '''0.07088'
$ws.Range("E9").Value = '  -3.22%  '
$ws.Range("D10").Value = '''0.8710'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").Value = '''0.07777'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = '''19.22'
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").Value = '1.822.37'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '''5.266'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").Value = '''6.322'
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("D16").Value = '''86.42'
$ws.Range("E16").Value = '  -5.89%  '
$ws.Range("D17").Value = '''1.010'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '''0.000008520'
$ws.Range("E18").Value = '  -5.06%  '
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = '26.490.63'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").Value = '''14.17'
$ws.Range("E21").Value = '  -4.28%  '
$ws.Range("D22").Value = '''4.953'
$ws.Range("E22").Value = '  -3.54%  '
$ws.Range("D23").Value = '2.038.44'
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").Value = '''10.35'
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").Value = '''1.977'
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("D26").Value = '''150.14'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").Value = '''17.82'
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").Value = '''1.991'
$ws.Range("E28").Value = '  -3.14%  '
$ws.Range("D29").Value = '''112.86'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").Value = '''4.864'
$ws.Range("E30").Value = '  -4.70%  '
$ws.Range("D31").Value = '''0.08641'
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("D32").Value = '''3.117'
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.435'
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7255'
$ws.Range("E34").Value = '  -5.42%  '
$ws.Range("D35").Value = '''1.111'
$ws.Range("E35").Value = '  -5.65%  '
$ws.Range("D36").Value = '''2.504'
$ws.Range("E36").Value = '  -7.69%  '
$ws.Range("D37").Value = '''1.077'
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").Value = '''0.01907'
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("D39").Value = '''0.05067'
$ws.Range("E39").Value = '  -3.45%  '
$ws.Range("D40").Value = '''2.868'
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("D41").Value = '''6.874'
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("D42").Value = '''0.4896'
$ws.Range("E42").Value = '  -4.96%  '
$ws.Range("D43").Value = '''0.1565'
$ws.Range("E43").Value = '  -5.01%  '
$ws.Range("D44").Value = '''8.111'
$ws.Range("E44").Value = '  -3.75%  '
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").Value = '''0.4587'
$ws.Range("E46").Value = '  -4.68%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''9.962'
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''101.06'
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("D49").Value = '''1.575'
$ws.Range("E49").Value = '  -4.52%  '
$ws.Range("D50").Value = '''0.05988'
$ws.Range("E50").Value = '  -3.75%  '
$ws.Range("D51").Value = '''63.26'
$ws.Range("E51").Value = '  -3.39%  '
